$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.479.57'
$ws.Range('E2').Value = '  +1.69%  '
$ws.Range('D3').Value = '1.845.65'
$ws.Range('E3').Value = '  +0.24%  '
$style = $ws.Range('D4').Style
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = $style
$ws.Range('E4').Value = '  +0.26%  '
$style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '260.41'
$ws.Range('D5').Style = $style
$ws.Range('E5').Value = '  -6.47%  '
$style = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('D6').Style = $style
$ws.Range('E6').Value = '  +0.27%  '
$style = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5216'
$ws.Range('D7').Style = $style
$ws.Range('E7').Value = '  +2.16%  '
$style = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3276'
$ws.Range('D8').Style = $style
$ws.Range('E8').Value = '  -6.34%  '
$style = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06748'
$ws.Range('D9').Style = $style
$ws.Range('E9').Value = '  -1.01%  '
$ws.Range('E10').Value = '  -2.54%  '
$style = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7748'
$ws.Range('D11').Style = $style
$ws.Range('E11').Value = '  -3.98%  '
$style = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07679'
$ws.Range('D12').Style = $style
$ws.Range('E12').Value = '  -1.22%  '
$ws.Range('D13').Value = '1.869.79'
$ws.Range('E13').Value = '  +1.76%  '
$style = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '88.43'
$ws.Range('D14').Style = $style
$ws.Range('E14').Value = '  +0.16%  '
$style = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.055'
$ws.Range('D15').Style = $style
$ws.Range('E15').Value = '  -0.52%  '
$ws.Range('E16').Value = '  +0.26%  '
$style = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.17'
$ws.Range('D17').Style = $style
$ws.Range('E17').Value = '  +0.03%  '
$style = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.002'
$ws.Range('D18').Style = $style
$ws.Range('E18').Value = '  +0.23%  '
$style = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007916'
$ws.Range('D19').Style = $style
$ws.Range('D20').Value = '26.502.72'
$ws.Range('E20').Value = '  +1.62%  '
$ws.Range('D21').Value = '2.098.66'
$ws.Range('E21').Value = '  +1.47%  '
$style = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.604'
$ws.Range('D22').Style = $style
$ws.Range('E22').Value = '  -3.64%  '
$style = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.615'
$ws.Range('D23').Style = $style
$ws.Range('E23').Value = '  -4.69%  '
$style = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.004'
$ws.Range('D24').Style = $style
$ws.Range('E24').Value = '  -3.33%  '
$style = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.338'
$ws.Range('D25').Style = $style
$ws.Range('E25').Value = '  -1.32%  '
$style = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '145.22'
$ws.Range('D26').Style = $style
$ws.Range('E26').Value = '  +1.09%  '
$style = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.640'
$ws.Range('D27').Style = $style
$ws.Range('E27').Value = '  -1.36%  '
$style = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.05'
$ws.Range('D28').Style = $style
$ws.Range('E28').Value = '  -0.85%  '
$style = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '111.54'
$ws.Range('D29').Style = $style
$ws.Range('E29').Value = '  +1.64%  '
$style = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.236'
$ws.Range('D30').Style = $style
$ws.Range('E30').Value = '  -2.88%  '
$style = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.196'
$ws.Range('D31').Style = $style
$ws.Range('E31').Value = '  -2.17%  '
$style = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08760'
$ws.Range('D32').Style = $style
$ws.Range('E32').Value = '  -0.55%  '
$style = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04855'
$ws.Range('D33').Style = $style
$ws.Range('E33').Value = '  -0.20%  '
$style = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.141'
$ws.Range('D34').Style = $style
$ws.Range('E34').Value = '  -1.99%  '
$style = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.869'
$ws.Range('D35').Style = $style
$ws.Range('E35').Value = '  +0.07%  '
$style = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7085'
$ws.Range('D36').Style = $style
$ws.Range('E36').Value = '  -2.87%  '
$style = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.104'
$ws.Range('D37').Style = $style
$ws.Range('E37').Value = '  -3.41%  '
$style = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01815'
$ws.Range('D38').Style = $style
$ws.Range('E38').Value = '  -1.99%  '
$style = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.232'
$ws.Range('D39').Style = $style
$ws.Range('E39').Value = '  -6.12%  '
$style = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.4952'
$ws.Range('D40').Style = $style
$ws.Range('E40').Value = '  -3.97%  '
$style = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '112.86'
$ws.Range('D41').Style = $style
$ws.Range('E41').Value = '  -3.73%  '
$ws.Range('E42').Value = '  -5.12%  '
$style = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.092'
$ws.Range('D43').Style = $style
$ws.Range('E43').Value = '  -2.72%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$style = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.002'
$ws.Range('D44').Style = $style
$ws.Range('E44').Value = '  +0.36%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$style = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.816'
$ws.Range('D45').Style = $style
$ws.Range('E45').Value = '  -2.47%  '
$style = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4293'
$ws.Range('D46').Style = $style
$ws.Range('E46').Value = '  -5.14%  '
$style = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1291'
$ws.Range('D47').Style = $style
$ws.Range('E47').Value = '  -5.28%  '
$style = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.159'
$ws.Range('D48').Style = $style
$ws.Range('E48').Value = '  -1.90%  '
$style = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05917'
$ws.Range('D49').Style = $style
$ws.Range('E49').Value = '  +0.03%  '
$style = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '35.39'
$ws.Range('D50').Style = $style
$ws.Range('E50').Value = '  -2.12%  '
$style = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.439'
$ws.Range('D51').Style = $style
$ws.Range('E51').Value = '  -3.66%  '
